$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.003.17'
$ws.Range("E2").Value = '  +1.36%  '
$ws.Range("D3").Value = '3.323.18'
$ws.Range("E3").Value = '  +6.26%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.79'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.99'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.47%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.322.20'
$ws.Range("E8").Value = '  +6.28%  '
$ws.Range("E9").Value = '  +1.48%  '
$ws.Range("E10").Value = '  +3.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.54'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +6.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.476'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.99%  '
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("E14").Value = '  +1.90%  '
$ws.Range("D15").Value = '3.868.42'
$ws.Range("E15").Value = '  +6.27%  '
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("D17").Value = '3.325.87'
$ws.Range("E17").Value = '  +6.31%  '
$ws.Range("D18").Value = '64.084.61'
$ws.Range("E19").Value = '  +3.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '484.31'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.33'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("E22").Value = '  +5.94%  '
$ws.Range("E23").Value = '  +3.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.76'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +5.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.98'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.88%  '
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("E27").Value = '  +2.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.30'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.29'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +4.32%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '29.82'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +11.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.17'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +5.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.106'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.95%  '
$ws.Range("E34").Value = '  +1.58%  '
$ws.Range("E35").Value = '  +2.89%  '
$ws.Range("E36").Value = '  +3.35%  '
$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0764'
$ws.Range("E37").Value = '  +8.12%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '53.39'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0402'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '434.61'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.92%  '
$ws.Range("D41").Value = '3.056.20'
$ws.Range("E41").Value = '  +5.66%  '
$ws.Range("E42").Value = '  +3.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.46'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.57%  '
$ws.Range("E44").Value = '  -2.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.268'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.81%  '
$ws.Range("E46").Value = '  +4.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.66'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '36.26'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +13.02%  '
$ws.Range("E50").Value = '  +2.75%  '
$ws.Range("E51").Value = '  +1.63%  '
